$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2076976666666666
$ws.Range("H2").Value = 0.6230929999999999
$ws.Range("I2").Value = 0.1089421893552267
$ws.Range("J2").Value = 0.1089421893552267
$ws.Range("M2").Value = 24.91851366666667
$ws.Range("N2").Value = 74.75554099999999
$ws.Range("O2").Value = 0.2924799159147552
$ws.Range("P2").Value = 0.2924799159147553
$ws.Range("Q2").Value = 5.175517145368111
$ws.Range("R2").Value = 46.57965430831299
$ws.Range("S2").Value = 0.03186340238218604
$ws.Range("T2").Value = 0.03186340238218605

$ws.Range("G3").Value = 0.2076976666666666
$ws.Range("H3").Value = 0.6230929999999999
$ws.Range("I3").Value = 0.1089421893552267
$ws.Range("J3").Value = 0.1089421893552267
$ws.Range("O3").Value = 0.4753125595076708
$ws.Range("P3").Value = 0.4753125595076708
$ws.Range("R3").Value = 75.69714536132298
$ws.Range("S3").Value = 0.05178159086080212
$ws.Range("T3").Value = 0.05178159086080211

$ws.Range("G4").Value = 0.2076976666666666
$ws.Range("H4").Value = 0.6230929999999999
$ws.Range("I4").Value = 0.1089421893552267
$ws.Range("J4").Value = 0.1089421893552267
$ws.Range("M4").Value = 19.78346566666667
$ws.Range("N4").Value = 59.350397
$ws.Range("O4").Value = 0.232207524577574
$ws.Range("P4").Value = 0.232207524577574
$ws.Range("Q4").Value = 4.108979657546778
$ws.Range("R4").Value = 36.980816917921
$ws.Range("S4").Value = 0.02529719611223852
$ws.Range("T4").Value = 0.02529719611223852

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.698796666666667
$ws.Range("H5").Value = 5.09639
$ws.Range("I5").Value = 0.8910578106447734
$ws.Range("J5").Value = 0.8910578106447733
$ws.Range("M5").Value = 24.91851366666667
$ws.Range("N5").Value = 74.75554099999999
$ws.Range("O5").Value = 0.2924799159147552
$ws.Range("P5").Value = 0.2924799159147553
$ws.Range("Q5").Value = 42.33148795522111
$ws.Range("R5").Value = 380.9833915969899
$ws.Range("S5").Value = 0.2606165135325692
$ws.Range("T5").Value = 0.2606165135325692

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.698796666666667
$ws.Range("H6").Value = 5.09639
$ws.Range("I6").Value = 0.8910578106447734
$ws.Range("J6").Value = 0.8910578106447733
$ws.Range("O6").Value = 0.4753125595076708
$ws.Range("P6").Value = 0.4753125595076708
$ws.Range("Q6").Value = 68.79340013769888
$ws.Range("R6").Value = 619.1406012392899
$ws.Range("S6").Value = 0.4235309686468687
$ws.Range("T6").Value = 0.4235309686468686

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.698796666666667
$ws.Range("H7").Value = 5.09639
$ws.Range("I7").Value = 0.8910578106447734
$ws.Range("J7").Value = 0.8910578106447733
$ws.Range("M7").Value = 19.78346566666667
$ws.Range("N7").Value = 59.350397
$ws.Range("O7").Value = 0.232207524577574
$ws.Range("P7").Value = 0.232207524577574
$ws.Range("Q7").Value = 33.60808552964778
$ws.Range("R7").Value = 302.47276976683
$ws.Range("S7").Value = 0.2069103284653355
$ws.Range("T7").Value = 0.2069103284653355

